# Weekly refresh of the "Mora" price sheet.
# The underlying data rows (2-11) got re-pulled for the week, which re-ordered
# the existing rows (each row keeps its original field values, just moves to
# a different row position) and also corrected the "Origen" of one record
# from "Provincia de Linares" to "Provincia de Curicó".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for each row (2-11), taken from the previous data for the week
# (rows are simply reshuffled), with the row-3 Origen correction applied.
$data = @{
    2  = @{ D = 44238; M = 300; N = 3600; O = 4000; P = 3800; R = "Provincia de Curicó"; S = 1900 }
    3  = @{ D = 44194; M = 120; N = 3000; O = 3000; P = 3000; R = "Provincia de Curicó"; S = 1500 }
    4  = @{ D = 44174; M = 200; N = 3200; O = 3200; P = 3200; R = "Provincia de Curicó"; S = 1600 }
    5  = @{ D = 44232; M = 200; N = 3000; O = 3000; P = 3000; R = "Provincia de Curicó"; S = 1500 }
    6  = @{ D = 44236; M = 300; N = 3600; O = 4000; P = 3800; R = "Provincia de Curicó"; S = 1900 }
    7  = @{ D = 44188; M = 150; N = 3000; O = 3400; P = 3240; R = "Provincia de Linares"; S = 1620 }
    8  = @{ D = 44237; M = 100; N = 3600; O = 4000; P = 3800; R = "Provincia de Curicó"; S = 1900 }
    9  = @{ D = 44208; M = 85;  N = 3000; O = 3000; P = 3000; R = "Provincia de Linares"; S = 1500 }
    10 = @{ D = 44168; M = 170; N = 8000; O = 8000; P = 8000; R = "Provincia de Linares"; S = 4000 }
    11 = @{ D = 44231; M = 150; N = 3400; O = 3400; P = 3400; R = "Provincia de Curicó"; S = 1700 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("N$row").Value = $vals.N
    $ws.Range("O$row").Value = $vals.O
    $ws.Range("P$row").Value = $vals.P
    $ws.Range("R$row").Value = $vals.R
    $ws.Range("S$row").Value = $vals.S
}
